# Refactor verb selection logic / update frontend interaction edit:
# Append two new verb rows (会う "to meet" and 言う "to say") to the
# conjugation table on Sheet1, following the existing Dictionary/Te/Ta/
# Nai/Masu/Volitional column layout, and move the on-screen view/selection
# to reflect where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the look & feel (style + row height) of the last existing data
#     row (87) down onto the two new rows before filling in values, so the
#     new rows render identically to the rest of the table.
$ws.Range("A87:F87").Copy()
$ws.Range("A88:F89").PasteSpecial(-4122)   # xlPasteFormats

# --- New row 88: 会う (au / "to meet") conjugations
$ws.Range("A88").Value = "会う"
$ws.Range("B88").Value = "会って"
$ws.Range("C88").Value = "あった"
$ws.Range("D88").Value = "会わない"
$ws.Range("E88").Value = "会います"
$ws.Range("F88").Value = "会おう"

# --- New row 89: 言う (iu / "to say") conjugations
$ws.Range("A89").Value = "言う"
$ws.Range("B89").Value = "言って"
$ws.Range("C89").Value = "言った"
$ws.Range("D89").Value = "言わない"
$ws.Range("E89").Value = "言います"
$ws.Range("F89").Value = "言おう"

# --- Match the explicit row heights used throughout the rest of the table.
$ws.Rows.Item(88).RowHeight = 18.75
$ws.Rows.Item(89).RowHeight = 18.75

# --- Update the saved view state: scrolled down a couple more rows, with
#     the active selection now on C85 instead of the old last cell F87.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 68
$ws.Range("C85").Select() | Out-Null
